# Ausarbeitung_Kapitelzuordnung.xlsx - Titel hinzugefügt, Quelle auf Seite 1 gefixt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix missing "Matrikel Nummer" (source) values for Andy Kruder's rows on sheet 1
$ws.Range("C29").Value = 7084984
$ws.Range("C30").Value = 7084984
$ws.Range("C31").Value = 7084984
$ws.Range("C32").Value = 7084984
$ws.Range("C52").Value = 7084984
$ws.Range("C53").Value = 7084984
$ws.Range("C54").Value = 7084984

# Update the view state to match the saved selection / scroll position
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("C70").Select()
